$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A:A").Insert()
$ws.Range("B1:L19").AutoFilter(8)
